$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price block (fecha serial 44448, variety "Lluteño") is inserted
# right before the existing block that starts at row 368, pushing all the
# rows from 368 downward by three rows (368-376 -> 371-379).
$ws.Rows("368:370").Insert()

# Columns that are identical across the three new rows.
$ws.Range("A368:A370").Value = 1
$ws.Range("B368:B370").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C368:C370").Value = "Arica y Parinacota"
$ws.Range("D368:D370").Value = 44448
$ws.Range("E368:E370").Value = 15
$ws.Range("F368:F370").Value = 100112024
$ws.Range("G368:G370").Value = "Choclo"
$ws.Range("H368:H370").Value = "Lluteño"
$ws.Range("O368:O370").Value = "Región de Arica y Parinacota"
$ws.Range("R368:R370").Value = "Hortaliza"

# Row 368: Primera
$ws.Range("I368").Value = "Primera"
$ws.Range("J368").Value = 70
$ws.Range("K368").Value = 32000
$ws.Range("L368").Value = 33000
$ws.Range("M368").Value = 32500
$ws.Range("N368").Value = "$/saco 50 unidades"
$ws.Range("P368").Value = 650
$ws.Range("Q368").Value = 50

# Row 369: Segunda
$ws.Range("I369").Value = "Segunda"
$ws.Range("J369").Value = 40
$ws.Range("K369").Value = 28000
$ws.Range("L369").Value = 30000
$ws.Range("M369").Value = 29000
$ws.Range("N369").Value = "$/saco 75 unidades"
$ws.Range("P369").Value = 387
$ws.Range("Q369").Value = 75

# Row 370: Tercera
$ws.Range("I370").Value = "Tercera"
$ws.Range("J370").Value = 50
$ws.Range("K370").Value = 24000
$ws.Range("L370").Value = 25000
$ws.Range("M370").Value = 24500
$ws.Range("N370").Value = "$/saco 100 unidades"
$ws.Range("P370").Value = 245
$ws.Range("Q370").Value = 100
